# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Sheets.Item("Error stats")
$ws2 = $wb.Sheets.Item("Residue stats")
$ws3 = $wb.Sheets.Item("Cal stats")

# --- Error stats: header relabel ---
$ws1.Range("R1").Value = "trans only noc"

# --- Error stats: data rows (cols B..U = 2..21) ---
$cols1 = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
$row1_4 = @(0.17789727827386054,0.61461270588568684,0.20999540722187857,0.52905714028872786,0.27147517449591774,0.74730970732864832,0.20999540722187857,0.52905714028872786,14.416661696809731,282.02243563285305,7.2273875425074117,144.99564823488595,14.271517668686256,283.4680085225985,7.2273875425074117,144.99564823488595,0.25327072596628042,0.57285531920565202,0.20074165402440952,0.38621836606251736)
for ($i = 0; $i -lt $cols1.Length; $i++) { $ws1.Cells.Item(4, $cols1[$i]).Value = $row1_4[$i] }
$row1_5 = @(0.42290273526674788,0.90324754175243749,0.38702259344750939,0.86863364675426313,0.50037234458763946,1.1020675756965639,0.38702259344750939,0.86863364675426313,14.425839437436602,286.07028333040728,7.9486026447623299,144.69170423656331,14.451512371099913,283.39170420874621,7.9486026447623299,144.69170423656331,0.53815219677763115,0.99928473208655955,0.35739244599688724,0.70269316327698628)
for ($i = 0; $i -lt $cols1.Length; $i++) { $ws1.Cells.Item(5, $cols1[$i]).Value = $row1_5[$i] }
$row1_6 = @(28.325383529852502,73.454045577923338,20.263774989255438,56.839424917636663,28.325383529852502,73.454045577923338,20.263774989255438,56.839424917636663,197.65522312805123,537.31567876127338,90.382233772846462,179.95013512759544,197.65522312805123,537.31567876127338,90.382233772846462,179.95013512759544,28.080008188075205,69.066140626855869,20.114019538485859,53.488912556240244)
for ($i = 0; $i -lt $cols1.Length; $i++) { $ws1.Cells.Item(6, $cols1[$i]).Value = $row1_6[$i] }
$row1_7 = @(50.238201276585791,261.94378445783059,35.137425843673363,156.7158296411321,50.238201276585791,261.94378445783059,35.137425843673363,156.7158296411321,153.47199000475959,575.01219048616701,67.98967504937049,179.86722798002737,153.47199000475959,575.01219048616701,67.98967504937049,179.86722798002737,49.670838073136842,192.27327476030777,35.054656482629809,148.14099229108331)
for ($i = 0; $i -lt $cols1.Length; $i++) { $ws1.Cells.Item(7, $cols1[$i]).Value = $row1_7[$i] }
$row1_8 = @(0.77021987371891576,2.6603888765214627,0.67589176773516291,1.5439673197988348,0.77459565135334385,2.7177749638241848,0.67589176773516291,1.5439673197988348,13.600682736914617,281.91401011271637,7.1659953340680493,144.47310143065357,13.51660379251588,282.55399317051325,7.1659953340680493,144.47310143065357,0.13186916238655375,0.39669922591158108,(7.7007609755504761 * [Math]::Pow(10, -2)),0.19615258837751368)
for ($i = 0; $i -lt $cols1.Length; $i++) { $ws1.Cells.Item(8, $cols1[$i]).Value = $row1_8[$i] }
$ws1.Cells.Item(8, 1).Value = "trans_only_hr_cal"

# --- Residue stats: header relabel ---
$ws2.Range("N1").Value = "trans only noc"

# --- Residue stats: data rows (cols B..P = 2..16) ---
$cols2 = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
$row2_3 = @((2.590066329459818 * [Math]::Pow(10, -7)),(1.6656841446641857 * [Math]::Pow(10, -6)),(3.4709830153805077 * [Math]::Pow(10, -6)),(2.590066329459818 * [Math]::Pow(10, -7)),(1.6656841446641857 * [Math]::Pow(10, -6)),(3.4709830153805077 * [Math]::Pow(10, -6)),(1.735526887163057 * [Math]::Pow(10, -6)),(1.4484609587310799 * [Math]::Pow(10, -5)),(1.8711875451724376 * [Math]::Pow(10, -2)),(1.735526887163057 * [Math]::Pow(10, -6)),(1.4484609587310799 * [Math]::Pow(10, -5)),(1.8711875451724376 * [Math]::Pow(10, -2)),(3.2709866412238952 * [Math]::Pow(10, -7)),(1.7054843955858453 * [Math]::Pow(10, -6)),(2.1353588105805418 * [Math]::Pow(10, -6)))
for ($i = 0; $i -lt $cols2.Length; $i++) { $ws2.Cells.Item(3, $cols2[$i]).Value = $row2_3[$i] }
$row2_4 = @((5.1073824295078025 * [Math]::Pow(10, -7)),(3.7822378109762738 * [Math]::Pow(10, -6)),(1.0238443056009985 * [Math]::Pow(10, -5)),(5.1073824295078025 * [Math]::Pow(10, -7)),(3.7822378109762738 * [Math]::Pow(10, -6)),(1.0238443056009985 * [Math]::Pow(10, -5)),(8.1754253385372231 * [Math]::Pow(10, -7)),(6.0890727413906567 * [Math]::Pow(10, -6)),(1.9104287874080948 * [Math]::Pow(10, -2)),(8.1754253385372231 * [Math]::Pow(10, -7)),(6.0890727413906567 * [Math]::Pow(10, -6)),(1.9104287874080948 * [Math]::Pow(10, -2)),(4.1110891999973994 * [Math]::Pow(10, -7)),(3.5513302502115253 * [Math]::Pow(10, -6)),(7.8427211582247019 * [Math]::Pow(10, -6)))
for ($i = 0; $i -lt $cols2.Length; $i++) { $ws2.Cells.Item(4, $cols2[$i]).Value = $row2_4[$i] }
$row2_5 = @((2.3230072688051068 * [Math]::Pow(10, -3)),(1.7169203218254957 * [Math]::Pow(10, -2)),(2.8052910813942691 * [Math]::Pow(10, -2)),(2.3230072688051068 * [Math]::Pow(10, -3)),(1.7169203218254957 * [Math]::Pow(10, -2)),(2.8052910813942691 * [Math]::Pow(10, -2)),(8.0745811880826076 * [Math]::Pow(10, -3)),(3.3399417000498234 * [Math]::Pow(10, -2)),(5.080602354241269 * [Math]::Pow(10, -2)),(8.0745811880826076 * [Math]::Pow(10, -3)),(3.3399417000498234 * [Math]::Pow(10, -2)),(5.080602354241269 * [Math]::Pow(10, -2)),(2.2109766895226008 * [Math]::Pow(10, -3)),(1.629398185977218 * [Math]::Pow(10, -2)),(2.7495250950826583 * [Math]::Pow(10, -2)))
for ($i = 0; $i -lt $cols2.Length; $i++) { $ws2.Cells.Item(5, $cols2[$i]).Value = $row2_5[$i] }
$row2_6 = @((2.1400510819570254 * [Math]::Pow(10, -3)),(2.0404362790901673 * [Math]::Pow(10, -2)),(4.9055255278758932 * [Math]::Pow(10, -2)),(2.1400510819570254 * [Math]::Pow(10, -3)),(2.0404362790901673 * [Math]::Pow(10, -2)),(4.9055255278758932 * [Math]::Pow(10, -2)),(6.307910416767838 * [Math]::Pow(10, -3)),(3.0439001022486936 * [Math]::Pow(10, -2)),(5.3981245480927448 * [Math]::Pow(10, -2)),(6.307910416767838 * [Math]::Pow(10, -3)),(3.0439001022486936 * [Math]::Pow(10, -2)),(5.3981245480927448 * [Math]::Pow(10, -2)),(2.1911950505628316 * [Math]::Pow(10, -3)),(1.6916280602721375 * [Math]::Pow(10, -2)),(4.0674871227342868 * [Math]::Pow(10, -2)))
for ($i = 0; $i -lt $cols2.Length; $i++) { $ws2.Cells.Item(6, $cols2[$i]).Value = $row2_6[$i] }
$row2_7 = @((1.2609555872881945 * [Math]::Pow(10, -6)),(6.6391451192545838 * [Math]::Pow(10, -6)),(2.4751039114775086 * [Math]::Pow(10, -5)),(1.2609555872881945 * [Math]::Pow(10, -6)),(6.6391451192545838 * [Math]::Pow(10, -6)),(2.4751039114775086 * [Math]::Pow(10, -5)),(2.1684066611749663 * [Math]::Pow(10, -6)),(1.1390501305779754 * [Math]::Pow(10, -5)),(1.911639392354585 * [Math]::Pow(10, -2)),(2.1684066611749663 * [Math]::Pow(10, -6)),(1.1390501305779754 * [Math]::Pow(10, -5)),(1.911639392354585 * [Math]::Pow(10, -2)),(4.6500833307231552 * [Math]::Pow(10, -8)),(2.8757943826056683 * [Math]::Pow(10, -7)),(4.8852248574876071 * [Math]::Pow(10, -7)))
for ($i = 0; $i -lt $cols2.Length; $i++) { $ws2.Cells.Item(7, $cols2[$i]).Value = $row2_7[$i] }
$ws2.Cells.Item(7, 1).Value = "trans_only_hr_cal"

# --- Cal stats: column B (RMS residue fraction) ---
$ws3.Cells.Item(2, 2).Value = (4.8777444372937878 * [Math]::Pow(10, -3))
$ws3.Cells.Item(3, 2).Value = (1.0815853169844849 * [Math]::Pow(10, -2))
$ws3.Cells.Item(4, 2).Value = 0.20810046205215704
$ws3.Cells.Item(5, 2).Value = 0.36649933201268786
$ws3.Cells.Item(6, 2).Value = (4.8777444372937644 * [Math]::Pow(10, -3))
$ws3.Cells.Item(7, 2).Value = (7.2722084314463488 * [Math]::Pow(10, -3))
$ws3.Cells.Item(8, 2).Value = (1.5690344591707652 * [Math]::Pow(10, -3))

# --- Cal stats: columns E/F (rows 6-8) ---
$ws3.Cells.Item(6, 5).Value = 0.31317047098370665
$ws3.Cells.Item(6, 6).Value = 0.19812485817603456
$ws3.Cells.Item(7, 5).Value = 2.282088115369139
$ws3.Cells.Item(7, 6).Value = 2.0774036731220922
$ws3.Cells.Item(8, 5).Value = 0.14180585410840837
$ws3.Cells.Item(8, 6).Value = (8.5407365701880067 * [Math]::Pow(10, -2))

# --- Cal stats: new row 8 (A/C/D) ---
$ws3.Cells.Item(8, 1).Value = "output/trans_only_hr_cal"
$ws3.Cells.Item(8, 3).Value = 30
$ws3.Cells.Item(8, 4).Value = 0

# --- Cal stats: selection ---
$ws3.Range("B2:F8").Select()
